$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.206.15"
$ws.Range("E2").Value = "  +1.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.134.78"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.99"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.79"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.462"
$ws.Range("E8").Value = "  +4.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.33"
$ws.Range("E9").Value = "  +2.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  -0.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.410"
$ws.Range("E11").Value = "  +3.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.667.79"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.137"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.69"
$ws.Range("E14").Value = "  +1.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000165"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.198.60"
$ws.Range("E16").Value = "  +0.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.125.79"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.02"
$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.69"
$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.16"
$ws.Range("E20").Value = "  +3.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "355.51"
$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.95"
$ws.Range("E23").Value = "  +0.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.506"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  -1.04%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0887"
$ws.Range("E27").Value = "  -3.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.30"
$ws.Range("E28").Value = "  -1.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.21"
$ws.Range("E29").Value = "  -0.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.87"
$ws.Range("E30").Value = "  -0.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.41"
$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.01"
$ws.Range("E32").Value = "  +1.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("E33").Value = "  -2.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.14"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.09"
$ws.Range("E35").Value = "  -1.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "25.93"
$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.26"
$ws.Range("E37").Value = "  -0.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.71"
$ws.Range("E38").Value = "  +5.39%  "

$ws.Range("E39").Value = "  +0.46%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.02"
$ws.Range("E40").Value = "  -3.54%  "

$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.700"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.75"
$ws.Range("E42").Value = "  +3.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.409.64"
$ws.Range("E43").Value = "  +3.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.175.19"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0269"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.983"
$ws.Range("E47").Value = "  +1.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.04"
$ws.Range("E48").Value = "  -0.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.92"
$ws.Range("E49").Value = "  -1.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.745"
$ws.Range("E50").Value = "  -2.00%  "

$ws.Range("E51").Value = "  +1.79%  "
